# "Generate Report for Archive"
#
# The localization run that produced this report moved forward: the three
# files that were sitting at "Ready for handoff" (81e72b09...md,
# 898fc443...yml, dc8e2780...yml) are now back "In Translation" on the
# Overview sheet as well as on each per-language detail sheet (zh-cn,
# de-de). Once no cell references the string "Ready for handoff" anymore,
# it naturally drops out of the workbook's shared-string table and every
# later shared-string index shifts down by one - exactly what the target
# OOXML shows.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet: rows 5-7, columns E (zh-cn) and F (de-de) ----
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($row in 5..7) {
    $wsOverview.Cells.Item($row, 5).Value = "In Translation"   # column E
    $wsOverview.Cells.Item($row, 6).Value = "In Translation"   # column F
}

# ---- zh-cn sheet: rows 5-7, column C (Status) ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($row in 5..7) {
    $wsZhCn.Cells.Item($row, 3).Value = "In Translation"       # column C
}

# ---- de-de sheet: rows 5-7, column C (Status) ----
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($row in 5..7) {
    $wsDeDe.Cells.Item($row, 3).Value = "In Translation"       # column C
}

# The now-shorter "In Translation" text re-shrinks the affected status
# columns relative to the old "Ready for handoff" width.
$wsOverview.Columns.Item(5).ColumnWidth = 13.04
$wsOverview.Columns.Item(6).ColumnWidth = 13.04
$wsZhCn.Columns.Item(3).ColumnWidth = 13.04
$wsDeDe.Columns.Item(3).ColumnWidth = 13.04
